$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap values between row 2 and row 4 for columns D, M, O, P, S
$ws.Range("D2").Value = 44672
$ws.Range("M2").Value = 8
$ws.Range("O2").Value = 180000
$ws.Range("P2").Value = 180000
$ws.Range("S2").Value = 180000

$ws.Range("D4").Value = 44993
$ws.Range("M4").Value = 14
$ws.Range("O4").Value = 200000
$ws.Range("P4").Value = 190000
$ws.Range("S4").Value = 190000
